$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# row 129 (hunk 0)
$ws.Range("H129").Value = 1196.5278
$ws.Range("J129").Value = 1287.8334
$ws.Range("L129").Value = 3863.5002
$ws.Range("N129").Value = -13863.5002

# row 132 (hunk 1)
$ws.Range("H132").Value = 1139.2709
$ws.Range("J132").Value = 2986.3333
$ws.Range("L132").Value = 8958.999899999999
$ws.Range("N132").Value = -14018.9999

# row 135 (hunk 2)
$ws.Range("H135").Value = 1287.25
$ws.Range("I135").Value = 1141.1
$ws.Range("J135").Value = 2018
$ws.Range("K135").Value = 10269.9
$ws.Range("L135").Value = 18162
$ws.Range("M135").Value = -7734.9
$ws.Range("N135").Value = -23232

# row 138 (hunk 3)
$ws.Range("H138").Value = 3975.12
$ws.Range("J138").Value = 4754.604
$ws.Range("L138").Value = 14263.812
$ws.Range("N138").Value = -24543.812

# row 140 (hunk 4)
$ws.Range("H140").Value = 146880
$ws.Range("J140").Value = 146880
$ws.Range("L140").Value = 146880
$ws.Range("N140").Value = -157240


$ws = $wb.Worksheets.Item("ARM")
# row 132 (hunk 5)
$ws.Range("H132").Value = 3371.0588
$ws.Range("I132").Value = 2942.5
$ws.Range("J132").Value = 4399.6
$ws.Range("K132").Value = 8827.5
$ws.Range("L132").Value = 13198.8
$ws.Range("M132").Value = -6297.5
$ws.Range("N132").Value = -18258.8

# row 133 (hunk 6)
$ws.Range("H133").Value = 42094.453
$ws.Range("J133").Value = 42094.453
$ws.Range("L133").Value = 42094.453
$ws.Range("N133").Value = -47154.453

# row 134 (hunk 7)
$ws.Range("H134").Value = 52592.332
$ws.Range("J134").Value = 52592.332
$ws.Range("L134").Value = 52592.332
$ws.Range("N134").Value = -62732.332

# row 135 (hunk 8)
$ws.Range("H135").Value = 48842.855
$ws.Range("J135").Value = 48842.855
$ws.Range("L135").Value = 48842.855
$ws.Range("N135").Value = -58982.855

# row 138 (hunk 9)
$ws.Range("H138").Value = 63050
$ws.Range("J138").Value = 63050
$ws.Range("L138").Value = 63050
$ws.Range("N138").Value = -73330

# row 141 (hunk 10)
$ws.Range("H141").Value = 51067.8
$ws.Range("J141").Value = 51067.8
$ws.Range("L141").Value = 51067.8
$ws.Range("N141").Value = -61427.8


$ws = $wb.Worksheets.Item("BSM")
# row 25 (hunk 11)
$ws.Range("H25").Value = 9921.333000000001
$ws.Range("I25").Value = 9921.333000000001
$ws.Range("K25").Value = 9921.333000000001
$ws.Range("M25").Value = -9686.333000000001

# row 99 (hunk 12)
$ws.Range("H99").Value = 1561.0667
$ws.Range("I99").Value = 1056
$ws.Range("K99").Value = 1056
$ws.Range("M99").Value = 442


$ws = $wb.Worksheets.Item("CRP")
# row 81 (hunk 13)
$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()

# row 84 (hunk 14)
$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()

# row 122 (hunk 15)
$ws.Range("H122").Value = 2072.2222
$ws.Range("I122").Value = 2364.4
$ws.Range("J122").Value = 1707
$ws.Range("K122").Value = 7093.200000000001
$ws.Range("L122").Value = 5121
$ws.Range("M122").Value = -4643.200000000001
$ws.Range("N122").Value = -10021

# row 133 (hunk 16)
$ws.Range("H133").Value = 49749.5
$ws.Range("J133").Value = 49749.5
$ws.Range("L133").Value = 49749.5
$ws.Range("N133").Value = -54809.5

# row 137 (hunk 17)
$ws.Range("H137").Value = 74590
$ws.Range("J137").Value = 74590
$ws.Range("L137").Value = 74590
$ws.Range("N137").Value = -84790

# row 138 (hunk 18)
$ws.Range("H138").Value = 74570
$ws.Range("J138").Value = 74570
$ws.Range("L138").Value = 74570
$ws.Range("N138").Value = -84850

# row 140 (hunk 19)
$ws.Range("H140").Value = 78500
$ws.Range("J140").Value = 78500
$ws.Range("L140").Value = 78500
$ws.Range("N140").Value = -88860

# row 141 (hunk 20)
$ws.Range("H141").Value = 39641
$ws.Range("J141").Value = 39621
$ws.Range("L141").Value = 39621
$ws.Range("N141").Value = -49981


$ws = $wb.Worksheets.Item("CUL")
# row 35 (hunk 21)
$ws.Range("H35").Value = 2049.75
$ws.Range("J35").Value = 2099.6667
$ws.Range("L35").Value = 6299.000100000001
$ws.Range("N35").Value = -6875.000100000001

# row 121 (hunk 22)
$ws.Range("H121").Value = 63748.75
$ws.Range("I121").Value = 500
$ws.Range("J121").Value = 72784.28999999999
$ws.Range("K121").Value = 1500
$ws.Range("L121").Value = 218352.87
$ws.Range("M121").Value = -190
$ws.Range("N121").Value = -220972.87


$ws = $wb.Worksheets.Item("GSM")
# row 132 (hunk 23)
$ws.Range("H132").Value = 1466.9032
$ws.Range("I132").Value = 1195.6154
$ws.Range("J132").Value = 2877.6
$ws.Range("K132").Value = 3586.8462
$ws.Range("L132").Value = 8632.799999999999
$ws.Range("M132").Value = -1056.8462
$ws.Range("N132").Value = -13692.8

# row 133 (hunk 24)
$ws.Range("H133").Value = 54033.332
$ws.Range("J133").Value = 54033.332
$ws.Range("L133").Value = 54033.332
$ws.Range("N133").Value = -64153.332

# row 135 (hunk 25)
$ws.Range("H135").Value = 47367.785
$ws.Range("J135").Value = 47367.785
$ws.Range("L135").Value = 47367.785
$ws.Range("N135").Value = -57507.785


$ws = $wb.Worksheets.Item("LTW")
# row 93 (hunk 26)
$ws.Range("H93").Value = 998
$ws.Range("I93").Value = 998
$ws.Range("K93").Value = 998
$ws.Range("M93").Value = 250

# row 133 (hunk 27)
$ws.Range("H133").Value = 61155.16
$ws.Range("J133").Value = 61155.16
$ws.Range("L133").Value = 61155.16
$ws.Range("N133").Value = -66215.16

# row 134 (hunk 28)
$ws.Range("H134").Value = 66550
$ws.Range("J134").Value = 66550
$ws.Range("L134").Value = 66550
$ws.Range("N134").Value = -76690

# row 138 (hunk 29)
$ws.Range("H138").Value = 58333.332
$ws.Range("J138").Value = 58333.332
$ws.Range("L138").Value = 58333.332
$ws.Range("N138").Value = -68613.33199999999

# row 140 (hunk 30)
$ws.Range("H140").Value = 83350
$ws.Range("J140").Value = 83350
$ws.Range("L140").Value = 83350
$ws.Range("N140").Value = -93710


$ws = $wb.Worksheets.Item("WVR")
# row 15 (hunk 31)
$ws.Range("H15").Value = 26405.6
$ws.Range("J15").Value = 26405.6
$ws.Range("L15").Value = 26405.6
$ws.Range("N15").Value = -26981.6

# row 122 (hunk 32)
$ws.Range("H122").Value = 50002260
$ws.Range("I122").Value = 125001000
$ws.Range("K122").Value = 375003000
$ws.Range("M122").Value = -375000550

# row 137 (hunk 33)
$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()

# row 138 (hunk 34)
$ws.Range("H138").Value = 109099
$ws.Range("J138").Value = 109099
$ws.Range("L138").Value = 109099
$ws.Range("N138").Value = -119379

# row 140 (hunk 35)
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()

